$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with new values (per diff)

# Row 2
$ws.Cells.Item(2, 1).Value = 45148.50694444445
$ws.Cells.Item(2, 2).Value = 15.854
$ws.Cells.Item(2, 3).Value = 10.494
$ws.Cells.Item(2, 4).Value = 3.735
$ws.Cells.Item(2, 5).Value = 33.762
$ws.Cells.Item(2, 6).Value = 26.097
$ws.Cells.Item(2, 7).Value = 12.381
$ws.Cells.Item(2, 8).Value = 37.784
$ws.Cells.Item(2, 9).Value = 19.197
$ws.Cells.Item(2, 10).Value = 7.826
$ws.Cells.Item(2, 11).Value = 11.635
$ws.Cells.Item(2, 12).Value = 13.329
$ws.Cells.Item(2, 13).Value = 13.955
$ws.Cells.Item(2, 14).Value = 3.981
$ws.Cells.Item(2, 15).Value = 12.407
$ws.Cells.Item(2, 16).Value = 17.155
$ws.Cells.Item(2, 17).Value = 10.858
$ws.Cells.Item(2, 18).Value = 3.133
$ws.Cells.Item(2, 19).Value = 2.01
$ws.Cells.Item(2, 20).Value = 181.475
$ws.Cells.Item(2, 21).Value = 34.412
$ws.Cells.Item(2, 22).Value = 11.452
$ws.Cells.Item(2, 23).Value = 22.298
$ws.Cells.Item(2, 24).Value = 11.285
$ws.Cells.Item(2, 25).Value = 3.17
$ws.Cells.Item(2, 26).Value = 19.349
$ws.Cells.Item(2, 27).Value = 10.115
$ws.Cells.Item(2, 28).Value = 9.106
$ws.Cells.Item(2, 29).Value = 10.918
$ws.Cells.Item(2, 30).Value = 14.036
$ws.Cells.Item(2, 31).Value = 3.313
$ws.Cells.Item(2, 32).Value = 33.905
$ws.Cells.Item(2, 33).Value = 6.108
$ws.Cells.Item(2, 34).Value = 14.317

# Row 3
$ws.Cells.Item(3, 1).Value = 45148.51388888889
$ws.Cells.Item(3, 2).Value = 6.726
$ws.Cells.Item(3, 3).Value = 4.368
$ws.Cells.Item(3, 4).Value = 1.452
$ws.Cells.Item(3, 5).Value = 14.491
$ws.Cells.Item(3, 6).Value = 10.941
$ws.Cells.Item(3, 7).Value = 5.21
$ws.Cells.Item(3, 8).Value = 22.653
$ws.Cells.Item(3, 9).Value = 8.144
$ws.Cells.Item(3, 10).Value = 3.334
$ws.Cells.Item(3, 11).Value = 4.709
$ws.Cells.Item(3, 12).Value = 5.753
$ws.Cells.Item(3, 13).Value = 6.013
$ws.Cells.Item(3, 14).Value = 1.695
$ws.Cells.Item(3, 15).Value = 5.263
$ws.Cells.Item(3, 16).Value = 7.28
$ws.Cells.Item(3, 17).Value = 4.834
$ws.Cells.Item(3, 18).Value = 1.396
$ws.Cells.Item(3, 19).Value = 0.783
$ws.Cells.Item(3, 20).Value = 72.813
$ws.Cells.Item(3, 21).Value = 14.886
$ws.Cells.Item(3, 22).Value = 4.858
$ws.Cells.Item(3, 23).Value = 9.532
$ws.Cells.Item(3, 24).Value = 4.837
$ws.Cells.Item(3, 25).Value = 1.399
$ws.Cells.Item(3, 26).Value = 10.833
$ws.Cells.Item(3, 27).Value = 4.291
$ws.Cells.Item(3, 28).Value = 3.987
$ws.Cells.Item(3, 29).Value = 4.735
$ws.Cells.Item(3, 30).Value = 6.044
$ws.Cells.Item(3, 31).Value = 1.246
$ws.Cells.Item(3, 32).Value = 21.172
$ws.Cells.Item(3, 33).Value = 2.512
$ws.Cells.Item(3, 34).Value = 6.075

# Row 4
$ws.Cells.Item(4, 1).Value = 45148.52083333334
$ws.Cells.Item(4, 2).Value = 21.139
$ws.Cells.Item(4, 3).Value = 15.429
$ws.Cells.Item(4, 4).Value = 1.466
$ws.Cells.Item(4, 5).Value = 45.926
$ws.Cells.Item(4, 6).Value = 37.234
$ws.Cells.Item(4, 7).Value = 16.562
$ws.Cells.Item(4, 8).Value = 60.976
$ws.Cells.Item(4, 9).Value = 25.596
$ws.Cells.Item(4, 10).Value = 11.268
$ws.Cells.Item(4, 11).Value = 16.631
$ws.Cells.Item(4, 12).Value = 18.408
$ws.Cells.Item(4, 13).Value = 19.442
$ws.Cells.Item(4, 14).Value = 5.313
$ws.Cells.Item(4, 15).Value = 16.542
$ws.Cells.Item(4, 16).Value = 23.451
$ws.Cells.Item(4, 17).Value = 14.059
$ws.Cells.Item(4, 18).Value = 1.052
$ws.Cells.Item(4, 19).Value = 0.968
$ws.Cells.Item(4, 20).Value = 244.465
$ws.Cells.Item(4, 21).Value = 46.18
$ws.Cells.Item(4, 22).Value = 15.269
$ws.Cells.Item(4, 23).Value = 30.894
$ws.Cells.Item(4, 24).Value = 16.182
$ws.Cells.Item(4, 25).Value = 2.63
$ws.Cells.Item(4, 26).Value = 30.142
$ws.Cells.Item(4, 27).Value = 13.487
$ws.Cells.Item(4, 28).Value = 12
$ws.Cells.Item(4, 29).Value = 14.136
$ws.Cells.Item(4, 30).Value = 19.343
$ws.Cells.Item(4, 31).Value = 0.766
$ws.Cells.Item(4, 32).Value = 55.414
$ws.Cells.Item(4, 33).Value = 8.523
$ws.Cells.Item(4, 34).Value = 19.09

# Row 5
$ws.Cells.Item(5, 1).Value = 45148.52777777778
$ws.Cells.Item(5, 2).Value = 16.33
$ws.Cells.Item(5, 3).Value = 11.93
$ws.Cells.Item(5, 4).Value = 1.1
$ws.Cells.Item(5, 5).Value = 35.5
$ws.Cells.Item(5, 6).Value = 28.78
$ws.Cells.Item(5, 7).Value = 12.79
$ws.Cells.Item(5, 8).Value = 51.18
$ws.Cells.Item(5, 9).Value = 19.78
$ws.Cells.Item(5, 10).Value = 8.73
$ws.Cells.Item(5, 11).Value = 12.85
$ws.Cells.Item(5, 12).Value = 14.24
$ws.Cells.Item(5, 13).Value = 15.05
$ws.Cells.Item(5, 14).Value = 4.11
$ws.Cells.Item(5, 15).Value = 12.78
$ws.Cells.Item(5, 16).Value = 18.15
$ws.Cells.Item(5, 17).Value = 10.87
$ws.Cells.Item(5, 18).Value = 0.79
$ws.Cells.Item(5, 19).Value = 0.72
$ws.Cells.Item(5, 20).Value = 187.25
$ws.Cells.Item(5, 21).Value = 35.79
$ws.Cells.Item(5, 22).Value = 11.8
$ws.Cells.Item(5, 23).Value = 23.95
$ws.Cells.Item(5, 24).Value = 12.53
$ws.Cells.Item(5, 25).Value = 2.02
$ws.Cells.Item(5, 26).Value = 24.68
$ws.Cells.Item(5, 27).Value = 10.42
$ws.Cells.Item(5, 28).Value = 9.279999999999999
$ws.Cells.Item(5, 29).Value = 10.92
$ws.Cells.Item(5, 30).Value = 14.96
$ws.Cells.Item(5, 31).Value = 0.55
$ws.Cells.Item(5, 32).Value = 46.69
$ws.Cells.Item(5, 33).Value = 6.58
$ws.Cells.Item(5, 34).Value = 14.75

# Delete row 6 (data no longer present; shrinks used range to A1:AH5)
$ws.Rows.Item(6).Delete()

# Widen columns 10, 22, 27, 34 from width 7 to width 8 characters
# (ColumnWidth uses "characters" units; 50/7 maps to a stored width of exactly 8)
$ws.Columns.Item(10).ColumnWidth = 50/7
$ws.Columns.Item(22).ColumnWidth = 50/7
$ws.Columns.Item(27).ColumnWidth = 50/7
$ws.Columns.Item(34).ColumnWidth = 50/7
